$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 35434
$ws.Range("I21").Value = 35607.6
$ws.Range("J21").Value = 35000
$ws.Range("K21").Value = 35607.6
$ws.Range("L21").Value = 35000
$ws.Range("M21").Value = -35139.6
$ws.Range("N21").Value = -35936
$ws.Range("H23").Value = 35434
$ws.Range("I23").Value = 35607.6
$ws.Range("J23").Value = 35000
$ws.Range("K23").Value = 35607.6
$ws.Range("L23").Value = 35000
$ws.Range("M23").Value = -35373.6
$ws.Range("N23").Value = -35468
$ws.Range("H32").Value = 22223120
$ws.Range("I32").Value = 47619336
$ws.Range("J32").Value = 1430
$ws.Range("K32").Value = 47619336
$ws.Range("L32").Value = 1430
$ws.Range("M32").Value = -47619010
$ws.Range("N32").Value = -2082
$ws.Range("H43").Value = 2054.3572
$ws.Range("I43").Value = 2281.5
$ws.Range("J43").Value = 1486.5
$ws.Range("K43").Value = 2281.5
$ws.Range("L43").Value = 1486.5
$ws.Range("M43").Value = -2212.5
$ws.Range("N43").Value = -1624.5
$ws.Range("I51").Value = 37038104
$ws.Range("J51").Value = 7511.1
$ws.Range("K51").Value = 37038104
$ws.Range("L51").Value = 7511.1
$ws.Range("M51").Value = -37037620
$ws.Range("N51").Value = -8479.1
$ws.Range("H64").Value = 2733.3333
$ws.Range("I64").Value = 2675
$ws.Range("J64").Value = 2850
$ws.Range("K64").Value = 2675
$ws.Range("L64").Value = 2850
$ws.Range("M64").Value = -2427
$ws.Range("N64").Value = -3346
$ws.Range("H67").Value = 2733.3333
$ws.Range("I67").Value = 2675
$ws.Range("J67").Value = 2850
$ws.Range("K67").Value = 2675
$ws.Range("L67").Value = 2850
$ws.Range("M67").Value = -1817
$ws.Range("N67").Value = -4566
$ws.Range("H132").Value = 12692.19
$ws.Range("I132").Value = 1964.3881
$ws.Range("J132").Value = 54972.35
$ws.Range("K132").Value = 5893.164299999999
$ws.Range("L132").Value = 164917.05
$ws.Range("M132").Value = -3363.164299999999
$ws.Range("N132").Value = -169977.05

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8633.671
$ws.Range("I32").Value = 7716.6304
$ws.Range("J32").Value = 19791
$ws.Range("K32").Value = 7716.6304
$ws.Range("L32").Value = 19791
$ws.Range("M32").Value = -7429.6304
$ws.Range("N32").Value = -20365
$ws.Range("H122").Value = 1798.1154
$ws.Range("I122").Value = 1789.591
$ws.Range("J122").Value = 1845
$ws.Range("K122").Value = 5368.772999999999
$ws.Range("L122").Value = 5535
$ws.Range("M122").Value = -2918.772999999999
$ws.Range("N122").Value = -10435
$ws.Range("H133").Value = 35729
$ws.Range("J133").Value = 35729
$ws.Range("L133").Value = 35729
$ws.Range("N133").Value = -40789

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 14126.533
$ws.Range("I12").Value = 1475
$ws.Range("J12").Value = 18727.092
$ws.Range("K12").Value = 1475
$ws.Range("L12").Value = 18727.092
$ws.Range("M12").Value = -1305
$ws.Range("N12").Value = -19067.092
$ws.Range("H132").Value = 22737.848
$ws.Range("I132").Value = 930.58
$ws.Range("K132").Value = 2791.74
$ws.Range("M132").Value = -261.7400000000002

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 148.81482
$ws.Range("I12").Value = 183.9
$ws.Range("J12").Value = 128.17647
$ws.Range("K12").Value = 551.7
$ws.Range("L12").Value = 384.52941
$ws.Range("M12").Value = -378.7
$ws.Range("N12").Value = -730.52941
$ws.Range("H54").Value = 3200
$ws.Range("J54").Value = 3200
$ws.Range("L54").Value = 9600
$ws.Range("N54").Value = -10718

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H33").Value = 19333.334
$ws.Range("J33").Value = 19333.334
$ws.Range("L33").Value = 19333.334
$ws.Range("N33").Value = -19837.334
$ws.Range("H102").Value = 2423.9092
$ws.Range("I102").Value = 2595.889
$ws.Range("J102").Value = 1650
$ws.Range("K102").Value = 2595.889
$ws.Range("L102").Value = 1650
$ws.Range("M102").Value = -973.8890000000001
$ws.Range("N102").Value = -4894
$ws.Range("H122").Value = 1049.1818
$ws.Range("I122").Value = 765
$ws.Range("J122").Value = 1211.5714
$ws.Range("K122").Value = 2295
$ws.Range("L122").Value = 3634.7142
$ws.Range("M122").Value = 155
$ws.Range("N122").Value = -8534.7142
$ws.Range("H126").Value = 22762.2
$ws.Range("I126").Value = 37104
$ws.Range("J126").Value = 1249.5
$ws.Range("K126").Value = 111312
$ws.Range("L126").Value = 3748.5
$ws.Range("M126").Value = -108842
$ws.Range("N126").Value = -8688.5
$ws.Range("H132").Value = 2527.5122
$ws.Range("I132").Value = 1725.4445
$ws.Range("K132").Value = 5176.333500000001
$ws.Range("M132").Value = -2646.333500000001

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2302.6943
$ws.Range("J7").Value = 2853.8572
$ws.Range("L7").Value = 2853.8572
$ws.Range("N7").Value = -3077.8572
$ws.Range("H22").Value = 438.26315
$ws.Range("I22").Value = 349
$ws.Range("J22").Value = 631.6667
$ws.Range("K22").Value = 349
$ws.Range("L22").Value = 631.6667
$ws.Range("M22").Value = -54
$ws.Range("N22").Value = -1221.6667
$ws.Range("H27").Value = 438.26315
$ws.Range("I27").Value = 349
$ws.Range("J27").Value = 631.6667
$ws.Range("K27").Value = 349
$ws.Range("L27").Value = 631.6667
$ws.Range("M27").Value = -242
$ws.Range("N27").Value = -845.6667
$ws.Range("H40").Value = 5138.2
$ws.Range("J40").Value = 8196.25
$ws.Range("L40").Value = 8196.25
$ws.Range("N40").Value = -8468.25
$ws.Range("H41").Value = 29000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 29000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 29000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -29876
$ws.Range("H46").Value = 5957
$ws.Range("I46").Value = 3030
$ws.Range("J46").Value = 8233.556
$ws.Range("K46").Value = 3030
$ws.Range("L46").Value = 8233.556
$ws.Range("M46").Value = -2842
$ws.Range("N46").Value = -8609.556
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H122").Value = 39859.04
$ws.Range("I122").Value = 51297
$ws.Range("J122").Value = 1732.5
$ws.Range("K122").Value = 153891
$ws.Range("L122").Value = 5197.5
$ws.Range("M122").Value = -151441
$ws.Range("N122").Value = -10097.5
$ws.Range("H126").Value = 2302.6943
$ws.Range("J126").Value = 2853.8572
$ws.Range("L126").Value = 8561.571599999999
$ws.Range("N126").Value = -13501.5716

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 757.6
$ws.Range("I113").Value = 756.6667
$ws.Range("J113").Value = 760.4
$ws.Range("K113").Value = 2270.0001
$ws.Range("L113").Value = 2281.2
$ws.Range("M113").Value = -100.0001000000002
$ws.Range("N113").Value = -6621.2
$ws.Range("H122").Value = 1905389.9
$ws.Range("I122").Value = 2198418.8
$ws.Range("J122").Value = 702
$ws.Range("K122").Value = 6595256.399999999
$ws.Range("L122").Value = 2106
$ws.Range("M122").Value = -6592806.399999999
$ws.Range("N122").Value = -7006
$ws.Range("H126").Value = 1472156
$ws.Range("I126").Value = 1962342.4
$ws.Range("K126").Value = 5887027.199999999
$ws.Range("M126").Value = -5884557.199999999
$ws.Range("H136").Value = 256693.86
$ws.Range("I136").Value = 306196.8
$ws.Range("J136").Value = 2107.2856
$ws.Range("K136").Value = 918590.3999999999
$ws.Range("L136").Value = 6321.8568
$ws.Range("M136").Value = -916040.3999999999
$ws.Range("N136").Value = -11421.8568
